$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for B2:O25 (columns B..O, rows 2..25), row-major order.
$data = @(
    @(16.82361958449795, 8.828868593120227, 5.987990474260265, 10.70594448902044, 0, 3.640681020734384, 0, 23.12352289791189, 0, 0, 9.993645275870515, 0, 17.45162601311849, 24.14055000092353),
    @(16.31107159179798, 8.473805340625189, 5.87005884353578, 10.7388079170951, 0, 3.643188961538602, 0, 23.2069061386077, 0, 0, 9.968754770430349, 0, 17.51015648793963, 24.15529775215094),
    @(15.99082084125332, 8.246311597960847, 5.798268398553287, 10.76055881878891, 0, 3.644810800605504, 0, 23.26370468744573, 0, 0, 9.955316409101874, 0, 17.54794695543648, 24.17111326079354),
    @(15.85912689255262, 8.151300566514884, 5.769214006851414, 10.76981804472281, 0, 3.64549238931397, 0, 23.28825416204098, 0, 0, 9.950307829929187, 0, 17.56381391452165, 24.17925441160765),
    @(15.83719366054018, 8.135387341894889, 5.7644029923355, 10.77137942603158, 0, 3.645606817395457, 0, 23.29241521578686, 0, 0, 9.949504509025212, 0, 17.56647685455767, 24.18070856131215),
    @(15.9890493041389, 8.24503946808322, 5.797875688299428, 10.76068209020327, 0, 3.644819908940905, 0, 23.26403009355483, 0, 0, 9.955246963334117, 0, 17.54815905018794, 24.17121619306687),
    @(16.64816221453172, 8.708455288103185, 5.947225753522576, 10.71694946617649, 0, 3.641528790481242, 0, 23.1511081263668, 0, 0, 9.984682630057714, 0, 17.47142355533343, 24.14422996410831),
    @(17.88843321288948, 9.539012183215007, 6.243077488957371, 10.64366315197704, 0, 3.635722064912869, 0, 22.97431899392965, 0, 0, 10.05685488903495, 0, 17.33558936605135, 24.14507463999816),
    @(18.75782408923492, 10.09825149268677, 6.459723977404036, 10.59741705790005, 0, 3.631846026189658, 0, 22.87193272868551, 0, 0, 10.11841391769996, 0, 17.24464227427628, 24.17857850719588),
    @(19.14244128208345, 10.34106178716267, 6.557620162825679, 10.57802687985333, 0, 3.63016650879909, 0, 22.83138352941395, 0, 0, 10.14820523769184, 0, 17.20517371085961, 24.20096089040041),
    @(19.28639148340697, 10.43130529723378, 6.594556730622096, 10.57092114708429, 0, 3.629542486411436, 0, 22.81689975143852, 0, 0, 10.15973728333318, 0, 17.19050057486697, 24.21046147977304),
    @(19.25546662249248, 10.41194602879372, 6.586608380683664, 10.57244095706799, 0, 3.629676349236659, 0, 22.81998026725575, 0, 0, 10.15724260522378, 0, 17.19364858218795, 24.20836982150674),
    @(19.15431885099341, 10.34852050509734, 6.560661875685374, 10.5774375388994, 0, 3.630114930478206, 0, 22.83017445151174, 0, 0, 10.1491489989402, 0, 17.20396108241705, 24.20172198558013),
    @(19.09213836049997, 10.30944770637774, 6.544750212174051, 10.5805289422021, 0, 3.630385131701199, 0, 22.83653229151241, 0, 0, 10.14422388755469, 0, 17.21031327767925, 24.197783381566),
    @(18.73245801466775, 10.08214675773372, 6.4533095296312, 10.59871740260653, 0, 3.631957465681261, 0, 22.87470439281838, 0, 0, 10.11650242410251, 0, 17.24725985943011, 24.17725931881033),
    @(18.50892220039756, 9.939708290036558, 6.39701608132823, 10.61029738147315, 0, 3.63294343713713, 0, 22.89966873505898, 0, 0, 10.09994968902006, 0, 17.2704122606565, 24.16649644367788),
    @(18.3793338871181, 9.856692290863414, 6.364577968058511, 10.61711293948273, 0, 3.633518424971118, 0, 22.91459456878484, 0, 0, 10.09059767471463, 0, 17.28390815536775, 24.16097829170503),
    @(18.33528720515493, 9.828398655319287, 6.353585997691426, 10.61944720174388, 0, 3.633714461843478, 0, 22.91974543800493, 0, 0, 10.08746040264145, 0, 17.28850844998246, 24.15922546370149),
    @(18.53282416943275, 9.954984117390811, 6.403015081523616, 10.60904862469819, 0, 3.632837663378145, 0, 22.89695252094297, 0, 0, 10.10169434305703, 0, 17.26792910405148, 24.16757259351618),
    @(19.18407541838501, 10.36719661478772, 6.568286954482846, 10.57596349200696, 0, 3.629985784056613, 0, 22.82715648862848, 0, 0, 10.15151953530805, 0, 17.20092465655912, 24.20364682358291),
    @(19.59976848788045, 10.62665907314277, 6.675500141881444, 10.55572132628413, 0, 3.628191682961508, 0, 22.78662153072454, 0, 0, 10.18554136832749, 0, 17.15872273857882, 24.23319535956576),
    @(19.37885473243634, 10.48909945596023, 6.618364455018947, 10.56639859590128, 0, 3.62914286560124, 0, 22.80778941033716, 0, 0, 10.16725204339045, 0, 17.18110159438564, 24.21687929142709),
    @(18.52202142611957, 9.948081421025876, 6.400303161794868, 10.60961269525644, 0, 3.632885458334301, 0, 22.89817873521656, 0, 0, 10.10090507248465, 0, 17.26905116193771, 24.16708398030737),
    @(17.55959085303909, 9.323069815358403, 6.162989621727928, 10.66215439510233, 0, 3.637224108391386, 0, 23.01733549434945, 0, 0, 10.03581025154235, 0, 17.37077643172135, 24.1390749061489)
)

$startRow = 2
$startCol = 2  # column B

for ($r = 0; $r -lt $data.Count; $r++) {
    $rowValues = $data[$r]
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $ws.Cells.Item($startRow + $r, $startCol + $c).Value = $rowValues[$c]
    }
}
